$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78; this shifts the existing rows 78-205
# down to 79-206 (matching Excel's native Insert behaviour), then we
# populate the newly-blank row 78 with the new data point.
$ws.Rows("78").Insert()

$ws.Range("A78").Value2 = 5
$ws.Range("B78").Value2 = "Macroferia Regional de Talca"
$ws.Range("C78").Value2 = "Maule"
$ws.Range("D78").Value2 = 44495
$ws.Range("E78").Value2 = 7
$ws.Range("F78").Value2 = 100114013
$ws.Range("G78").Value2 = "Zanahoria"
$ws.Range("H78").Value2 = "Sin especificar"
$ws.Range("I78").Value2 = "Primera"
$ws.Range("J78").Value2 = 400
$ws.Range("K78").Value2 = 8000
$ws.Range("L78").Value2 = 8000
$ws.Range("M78").Value2 = 8000
$ws.Range("N78").Value2 = "$/saco 20 kilos"
$ws.Range("O78").Value2 = "Región de Ñuble"
$ws.Range("P78").Value2 = 400
$ws.Range("Q78").Value2 = 20
$ws.Range("R78").Value2 = "Hortaliza"
